$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 119; all rows from 119 downward shift down by one.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new record's data.
$ws.Range("A119").Value = 4
$ws.Range("B119").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C119").Value = "Los Lagos"
$ws.Range("D119").Value = 44468
$ws.Range("E119").Value = 10
$ws.Range("F119").Value = 100112037
$ws.Range("G119").Value = "Cebollín"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 70
$ws.Range("K119").Value = 6000
$ws.Range("L119").Value = 6000
$ws.Range("M119").Value = 6000
$ws.Range("N119").Value = "$/paquete 36 unidades"
$ws.Range("O119").Value = "Región Metropolitana"
$ws.Range("P119").Value = 167
$ws.Range("Q119").Value = 36
$ws.Range("R119").Value = "Hortaliza"
